# Cotações atualizadas - 2025-10-03
# Append a new row (29) with the quotes for 2025-10-03, continuing the
# existing table in Sheet1 (dates/quotes in columns A-E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. the date number format) from the last existing
# row's date cell (A28) onto the new date cell (A29) before setting values.
$ws.Cells.Item(28, 1).Copy($ws.Cells.Item(29, 1))

$ws.Cells.Item(29, 1).Value = 45933
$ws.Cells.Item(29, 2).Value = "21,4463"
$ws.Cells.Item(29, 3).Value = "15,2675"
$ws.Cells.Item(29, 4).Value = "15,4193"
$ws.Cells.Item(29, 5).Value = "15,4193"
